# Edit script: update product row to use numeric ids instead of
# string codes (P_00001 -> 1, R_00001 -> 1), hide the helper JSON
# columns (I:P) while keeping column H visible, and reset the
# selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values -------------------------------------------------
# idPlato: was text "P_00001", now numeric 1
$ws.Range("A2").Value = 1

# idRestaurante: was text "R_00001", now numeric 1
$ws.Range("G2").Value = 1

# --- Column visibility ---------------------------------------------------
# Column H (8) stays visible; columns I:P (9-16) become hidden.
$ws.Range("I1:P1").EntireColumn.Hidden = $true

# --- View / selection ------------------------------------------------
$ws.Range("G2").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
